$wb = $excel.ActiveWorkbook
$dataSheet = $wb.Worksheets.Item("data")

# Insert the new "metadata" sheet right after "data"
$ws = $wb.Worksheets.Add($null, $dataSheet)
$ws.Name = "metadata"

# Match the page margins used on the rest of the workbook
$ws.PageSetup.LeftMargin = 54
$ws.PageSetup.RightMargin = 54
$ws.PageSetup.TopMargin = 72
$ws.PageSetup.BottomMargin = 72
$ws.PageSetup.HeaderMargin = 36
$ws.PageSetup.FooterMargin = 36

# ---- Header row (row 1) ----
$ws.Range("B1").Value = "data_name"
$ws.Range("C1").Value = "data_id"
$ws.Range("D1").Value = "data_version"
$ws.Range("E1").Value = "data_version_created"
$ws.Range("F1").Value = "panel_query_time"
$ws.Range("G1").Value = "panel_get_request"

# Reuse the same header style already used on the "data" sheet (bold, centered, bordered)
$dataSheet.Range("B1").Copy()
$ws.Range("B1:G1").PasteSpecial(-4122)

# ---- Data row (row 2) ----
$ws.Range("A2").Value = 0
$dataSheet.Range("A2").Copy()
$ws.Range("A2").PasteSpecial(-4122)

$ws.Range("B2").Value = "Polycystic liver disease"
$ws.Range("C2").Value = 3274
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "0.27"
$ws.Range("E2").Value = "2020-12-06T20:59:57.024826Z"
$ws.Range("F2").Value = "2021-10-05 14:35:15.886917"
$ws.Range("G2").Value = "https://panelapp.agha.umccr.org/api/v1/panels/3274/?format=json"

# ---- Refresh the F column "time_taken" timestamps on the "data" sheet ----
$dataSheet.Range("F2").Value = "2021-10-05 14:35:15.890096"
$dataSheet.Range("F3").Value = "2021-10-05 14:35:15.890104"
$dataSheet.Range("F4").Value = "2021-10-05 14:35:15.890107"
$dataSheet.Range("F5").Value = "2021-10-05 14:35:15.890109"
$dataSheet.Range("F6").Value = "2021-10-05 14:35:15.890113"
$dataSheet.Range("F7").Value = "2021-10-05 14:35:15.890115"
$dataSheet.Range("F8").Value = "2021-10-05 14:35:15.890118"
$dataSheet.Range("F9").Value = "2021-10-05 14:35:15.890120"
$dataSheet.Range("F10").Value = "2021-10-05 14:35:15.890123"
$dataSheet.Range("F11").Value = "2021-10-05 14:35:15.890126"
$dataSheet.Range("F12").Value = "2021-10-05 14:35:15.890128"

$ws.Range("A1").Select()
$dataSheet.Activate()
